# Periodic GitHub Actions refresh of the cryptocurrency tracking sheet:
# updates each coin's Price (column D) and 1h Volume % (column E) with the
# latest scraped figures. A handful of rows were also re-ranked relative
# to one another (ImmutableX / PancakeSwap / NEARProtocol, and
# VeChain / Maker), so their Coin name / Link / Price / Volume cells moved
# down a row.
#
# Price values are stored as plain text in this sheet (several rows use
# '.' as a thousands separator, e.g. "64.125.88"), so any new price that
# Excel would otherwise auto-recognise as a number is briefly forced to
# Text format before being written, then restored to the sheet's normal
# (unformatted) cell style so no visible formatting changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Reference cell with the sheet's normal/default style (never edited by
# this script) used to restore styling after temporarily forcing Text
# format on numeric-looking price cells.
$styleRef = $ws.Range('D4')

$ws.Range('D2').Value = '64.125.88'
$ws.Range('E2').Value = '  +0.68%  '
$ws.Range('D3').Value = '3.149.56'
$ws.Range('E3').Value = '  +0.77%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.18'
$ws.Range('D5').Style = $styleRef.Style
$ws.Range('E5').Value = '  +0.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.55'
$ws.Range('D6').Style = $styleRef.Style
$ws.Range('E6').Value = '  -0.30%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '3.143.76'
$ws.Range('E8').Value = '  +0.71%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.530'
$ws.Range('D9').Style = $styleRef.Style
$ws.Range('E9').Value = '  -0.46%  '
$ws.Range('E10').Value = '  +0.54%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.88'
$ws.Range('D11').Style = $styleRef.Style
$ws.Range('E11').Value = '  +3.29%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.461'
$ws.Range('D12').Style = $styleRef.Style
$ws.Range('E12').Value = '  -1.81%  '
$ws.Range('E13').Value = '  -1.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.21'
$ws.Range('D14').Style = $styleRef.Style
$ws.Range('E14').Value = '  +0.38%  '
$ws.Range('D15').Value = '3.671.24'
$ws.Range('E15').Value = '  +0.76%  '
$ws.Range('E16').Value = '  -1.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.27'
$ws.Range('D17').Style = $styleRef.Style
$ws.Range('E17').Value = '  +1.42%  '
$ws.Range('D18').Value = '63.948.69'
$ws.Range('E18').Value = '  +0.57%  '
$ws.Range('D19').Value = '3.145.48'
$ws.Range('E19').Value = '  +0.60%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '468.36'
$ws.Range('D20').Style = $styleRef.Style
$ws.Range('E20').Value = '  +0.66%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.38'
$ws.Range('D21').Style = $styleRef.Style
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.733'
$ws.Range('D22').Style = $styleRef.Style
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('E23').Value = '  -0.27%  '
$ws.Range('E24').Value = '  -1.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '81.39'
$ws.Range('D25').Style = $styleRef.Style
$ws.Range('E25').Value = '  -1.07%  '
$ws.Range('E26').Value = '  +6.02%  '
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.71'
$ws.Range('D28').Style = $styleRef.Style
$ws.Range('E28').Value = '  +7.97%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.71'
$ws.Range('D29').Style = $styleRef.Style
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.39'
$ws.Range('D30').Style = $styleRef.Style
$ws.Range('E30').Value = '  +7.46%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.24'
$ws.Range('D31').Style = $styleRef.Style
$ws.Range('E31').Value = '  +0.30%  '
$ws.Range('E32').Value = '  +0.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.71'
$ws.Range('D33').Style = $styleRef.Style
$ws.Range('E33').Value = '  +2.00%  '
$ws.Range('E34').Value = '  +2.23%  '
$ws.Range('D35').Value = '0.0₃0837'
$ws.Range('E35').Value = '  -4.57%  '
$ws.Range('E36').Value = '  +1.51%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.17'
$ws.Range('D37').Style = $styleRef.Style
$ws.Range('E37').Value = '  +0.74%  '
$ws.Range('E38').Value = '  -2.80%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.25'
$ws.Range('D39').Style = $styleRef.Style
$ws.Range('E39').Value = '  -5.28%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '461.83'
$ws.Range('D40').Style = $styleRef.Style
$ws.Range('E40').Value = '  +1.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '51.40'
$ws.Range('D41').Style = $styleRef.Style
$ws.Range('E41').Value = '  +0.80%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.20'
$ws.Range('D42').Style = $styleRef.Style
$ws.Range('E42').Value = '  +5.30%  '
$ws.Range('E43').Value = '  +5.22%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.929.92'
$ws.Range('E44').Value = '  +0.81%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0373'
$ws.Range('D45').Style = $styleRef.Style
$ws.Range('E45').Value = '  -0.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.62'
$ws.Range('D46').Style = $styleRef.Style
$ws.Range('E46').Value = '  +11.10%  '
$ws.Range('E47').Value = '  -2.31%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '130.19'
$ws.Range('D48').Style = $styleRef.Style
$ws.Range('E48').Value = '  +1.97%  '
$ws.Range('E49').Value = '  -0.04%  '
$ws.Range('E50').Value = '  +2.61%  '
$ws.Range('E51').Value = '  -0.80%  '
